$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price-only (column D) updates ---
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "269.53"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "22.89"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "6.389"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.06238"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "3.646"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "6.693"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.376"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.8359"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.01379"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08418"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.03418"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.03111"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.04687"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.006914"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1169"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.003206"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.01117"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00006281"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00000000751"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.9006"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.07388"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.00002202"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.01241"

# --- Full row updates (B, C, D, E) for rows 15-26 (ProBitToken moved to top) ---
$ws.Cells.Item(15, 2).Value = "ProBitToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.1264"
$ws.Cells.Item(15, 5).Value = "14ProBitTokenPROB"
$ws.Cells.Item(16, 2).Value = "BitMartToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.09330"
$ws.Cells.Item(16, 5).Value = "15BitMartTokenBMX"
$ws.Cells.Item(17, 2).Value = "MCDex"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.888"
$ws.Cells.Item(17, 5).Value = "16MCDexMCB"
$ws.Cells.Item(18, 2).Value = "BitForexToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.001717"
$ws.Cells.Item(18, 5).Value = "17BitForexTokenBF"
$ws.Cells.Item(19, 2).Value = "CoinExToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.04819"
$ws.Cells.Item(19, 5).Value = "18CoinExTokenCET"
$ws.Cells.Item(20, 2).Value = "TigerCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.006301"
$ws.Cells.Item(20, 5).Value = "19TigerCashTCH"
$ws.Cells.Item(21, 2).Value = "BitKan"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.001088"
$ws.Cells.Item(21, 5).Value = "20BitKanKAN"
$ws.Cells.Item(22, 2).Value = "HotbitToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.003325"
$ws.Cells.Item(22, 5).Value = "21HotbitTokenHTBWorstin24h"
$ws.Cells.Item(23, 2).Value = "NitroEx"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.0001501"
$ws.Cells.Item(23, 5).Value = "22NitroExNTX"
$ws.Cells.Item(24, 2).Value = "LEO"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "3.733"
$ws.Cells.Item(24, 5).Value = "23LEOLEO"
$ws.Cells.Item(25, 2).Value = "BTSEToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.370"
$ws.Cells.Item(25, 5).Value = "24BTSETokenBTSE"
$ws.Cells.Item(26, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.3404"
$ws.Cells.Item(26, 5).Value = "25BitpandaEcosystemTokenBEST"

Write-Output "done"